# Generate Report for Handback
#
# For the "9b517dad-f4cf-4ba8-b9ed-fa981a5cdcfe" handback row (row 7) on both
# the zh-cn and de-de sheets, a handback file has now actually been picked
# up, but it turned out to be stale (based on an older commit than the
# latest one available). Reflect that in the report:
#   - "Latest Target File"    (I7): fill in + hyperlink to the source .md
#   - "Latest Handback File"  (J7): the handback .xlf filename
#   - "Latest Handback DateTime" (K7): the handback timestamp
#   - "Error Detail"          (P7): explanation that the handback is stale

$wb = $excel.ActiveWorkbook

$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91f20638a67be4233f5a6d74eccd2721fbd26077/e2e/9b517dad-f4cf-4ba8-b9ed-fa981a5cdcfe.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/905cb14f5ec8a635b6e10edbfd15640e6c6c055d/e2e/9b517dad-f4cf-4ba8-b9ed-fa981a5cdcfe.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91f20638a67be4233f5a6d74eccd2721fbd26077/e2e/9b517dad-f4cf-4ba8-b9ed-fa981a5cdcfe.md."

$sheetsInfo = @(
    @{ Name = "zh-cn"; Handback = "9b517dad-f4cf-4ba8-b9ed-fa981a5cdcfe.d8282444ce11375edb0d75c44c0ab73295aac5e7.zh-cn.xlf"; DateTime = "2016-08-26 15:08:14" },
    @{ Name = "de-de"; Handback = "9b517dad-f4cf-4ba8-b9ed-fa981a5cdcfe.d8282444ce11375edb0d75c44c0ab73295aac5e7.de-de.xlf"; DateTime = "2016-08-26 15:08:22" }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Latest Target File: now resolved to the source markdown, linked back
    # to GitHub like the other "I" column entries on this sheet.
    $ws.Range("I7").Value = "9b517dad-f4cf-4ba8-b9ed-fa981a5cdcfe.md"
    $ws.Range("I7").Style = "HyperLink"
    $ws.Hyperlinks.Add($ws.Range("I7"), $targetMdUrl, "", "", "9b517dad-f4cf-4ba8-b9ed-fa981a5cdcfe.md")

    # Latest Handback File / DateTime
    $ws.Range("J7").Value = $info.Handback
    $ws.Range("K7").Value = $info.DateTime

    # Error Detail: handback exists, but it's against a stale commit.
    $ws.Range("P7").Value = $errorDetail
}
